$wb = $excel.ActiveWorkbook

# This script applies the updated market-price / profit figures produced by
# the scheduled Sagittarius Profits data-refresh run, sheet by sheet.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5394
$ws.Range("J2").Value = 5693
$ws.Range("L2").Value = 5693
$ws.Range("N2").Value = -5919
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 3500
$ws.Range("K4").Value = 3500
$ws.Range("M4").Value = -3386
$ws.Range("H15").Value = 1551.9803
$ws.Range("I15").Value = 1551.9803
$ws.Range("K15").Value = 4655.9409
$ws.Range("M15").Value = -4486.9409
$ws.Range("H19").Value = 1483.6
$ws.Range("I19").Value = 1677.6
$ws.Range("J19").Value = 1289.6
$ws.Range("K19").Value = 1677.6
$ws.Range("L19").Value = 1289.6
$ws.Range("M19").Value = -1502.6
$ws.Range("N19").Value = -1639.6
$ws.Range("H28").Value = 1601.4445
$ws.Range("I28").Value = 1859
$ws.Range("K28").Value = 1859
$ws.Range("M28").Value = -1374
$ws.Range("H49").Value = 875
$ws.Range("J49").Value = 950
$ws.Range("L49").Value = 2850
$ws.Range("N49").Value = -3122
$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4114
$ws.Range("H62").Value = 7999.3335
$ws.Range("I62").Value = 7998
$ws.Range("K62").Value = 7998
$ws.Range("M62").Value = -7374
$ws.Range("H65").Value = 7999.3335
$ws.Range("I65").Value = 7998
$ws.Range("K65").Value = 39990
$ws.Range("M65").Value = -36870
$ws.Range("H86").Value = 5727.9565
$ws.Range("I86").Value = 4999.2666
$ws.Range("K86").Value = 4999.2666
$ws.Range("M86").Value = -3876.2666
$ws.Range("H89").Value = 5727.9565
$ws.Range("I89").Value = 4999.2666
$ws.Range("K89").Value = 24996.333
$ws.Range("M89").Value = -19380.333
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H112").Value = 900.8889
$ws.Range("J112").Value = 900.8889
$ws.Range("L112").Value = 2702.6667
$ws.Range("N112").Value = -4918.6667
$ws.Range("H133").Value = 81499
$ws.Range("J133").Value = 81499
$ws.Range("L133").Value = 81499
$ws.Range("N133").Value = -91619
$ws.Range("H137").Value = 1339
$ws.Range("I137").Value = 733
$ws.Range("J137").Value = 2248
$ws.Range("K137").Value = 2199
$ws.Range("L137").Value = 6744
$ws.Range("M137").Value = 351
$ws.Range("N137").Value = -11844
$ws.Range("H138").Value = 3568.0454
$ws.Range("J138").Value = 3789.25
$ws.Range("L138").Value = 11367.75
$ws.Range("N138").Value = -21647.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 597.25
$ws.Range("I26").Value = 597.25
$ws.Range("K26").Value = 597.25
$ws.Range("M26").Value = -267.25
$ws.Range("H110").Value = 1258.1666
$ws.Range("I110").Value = 1258.1666
$ws.Range("K110").Value = 1258.1666
$ws.Range("M110").Value = 786.8334
$ws.Range("H113").Value = 142848.5
$ws.Range("J113").Value = 142848.5
$ws.Range("L113").Value = 142848.5
$ws.Range("N113").Value = -151526.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1245.25
$ws.Range("I23").Value = 2001
$ws.Range("J23").Value = 993.3333
$ws.Range("K23").Value = 2001
$ws.Range("L23").Value = 993.3333
$ws.Range("M23").Value = -1761
$ws.Range("N23").Value = -1473.3333
$ws.Range("H27").Value = 1245.25
$ws.Range("I27").Value = 2001
$ws.Range("J27").Value = 993.3333
$ws.Range("K27").Value = 2001
$ws.Range("L27").Value = 993.3333
$ws.Range("M27").Value = -1809
$ws.Range("N27").Value = -1377.3333
$ws.Range("H58").Value = 5014
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 5014
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 5014
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -5420
$ws.Range("H134").Value = 2476.5789
$ws.Range("I134").Value = 2417.4
$ws.Range("K134").Value = 7252.200000000001
$ws.Range("M134").Value = -4717.200000000001
$ws.Range("H136").Value = 5014
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5014
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 15042
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20142

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 155.85715
$ws.Range("I12").Value = 187.16667
$ws.Range("J12").Value = 132.375
$ws.Range("K12").Value = 561.50001
$ws.Range("L12").Value = 397.125
$ws.Range("M12").Value = -388.50001
$ws.Range("N12").Value = -743.125
$ws.Range("H23").Value = 16667047
$ws.Range("I23").Value = 20000370
$ws.Range("J23").Value = 433
$ws.Range("K23").Value = 60001110
$ws.Range("L23").Value = 1299
$ws.Range("M23").Value = -60000875
$ws.Range("N23").Value = -1769
$ws.Range("H33").Value = 100
$ws.Range("I33").Value = 100
$ws.Range("K33").Value = 600
$ws.Range("M33").Value = -317
$ws.Range("H37").Value = 43999.5
$ws.Range("J37").Value = 43999.5
$ws.Range("L37").Value = 131998.5
$ws.Range("N37").Value = -132222.5
$ws.Range("H40").Value = 115.888885
$ws.Range("I40").Value = 55.25
$ws.Range("K40").Value = 221
$ws.Range("M40").Value = -152
$ws.Range("H59").Value = 1947.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1947.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5842.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -6922.5
$ws.Range("H60").Value = 1032.6666
$ws.Range("I60").Value = 602
$ws.Range("K60").Value = 1806
$ws.Range("M60").Value = -1555
$ws.Range("H113").Value = 1226.0526
$ws.Range("J113").Value = 1216.8667
$ws.Range("L113").Value = 3650.6001
$ws.Range("N113").Value = -7990.6001
$ws.Range("H118").Value = 3302.2
$ws.Range("I118").Value = 3627.75
$ws.Range("K118").Value = 10883.25
$ws.Range("M118").Value = -9640.25
$ws.Range("H122").Value = 652.2857
$ws.Range("J122").Value = 613.4
$ws.Range("L122").Value = 5520.599999999999
$ws.Range("N122").Value = -10420.6
$ws.Range("H132").Value = 6453.6924
$ws.Range("I132").Value = 2488.889
$ws.Range("J132").Value = 15374.5
$ws.Range("K132").Value = 22400.001
$ws.Range("L132").Value = 138370.5
$ws.Range("M132").Value = -19870.001
$ws.Range("N132").Value = -143430.5
$ws.Range("H134").Value = 16142
$ws.Range("J134").Value = 18976
$ws.Range("L134").Value = 56928
$ws.Range("N134").Value = -67068

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 56450
$ws.Range("J95").Value = 56450
$ws.Range("L95").Value = 56450
$ws.Range("N95").Value = -61942
$ws.Range("H118").Value = 86211.5
$ws.Range("J118").Value = 86211.5
$ws.Range("L118").Value = 86211.5
$ws.Range("N118").Value = -89525.5
$ws.Range("H132").Value = 2541.0833
$ws.Range("I132").Value = 2421.6667
$ws.Range("J132").Value = 2899.3333
$ws.Range("K132").Value = 7265.000100000001
$ws.Range("L132").Value = 8697.999899999999
$ws.Range("M132").Value = -4735.000100000001
$ws.Range("N132").Value = -13757.9999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 6000
$ws.Range("I25").Value = 6000
$ws.Range("K25").Value = 6000
$ws.Range("M25").Value = -5770
$ws.Range("H40").Value = 4002.1667
$ws.Range("J40").Value = 4434.625
$ws.Range("L40").Value = 4434.625
$ws.Range("N40").Value = -4706.625
$ws.Range("H61").Value = 5649.7
$ws.Range("I61").Value = 4973.875
$ws.Range("J61").Value = 8353
$ws.Range("K61").Value = 4973.875
$ws.Range("L61").Value = 8353
$ws.Range("M61").Value = -4771.875
$ws.Range("N61").Value = -8757
$ws.Range("H113").Value = 5649.7
$ws.Range("I113").Value = 4973.875
$ws.Range("J113").Value = 8353
$ws.Range("K113").Value = 4973.875
$ws.Range("L113").Value = 8353
$ws.Range("M113").Value = -2803.875
$ws.Range("N113").Value = -12693
$ws.Range("H122").Value = 5963.9565
$ws.Range("I122").Value = 5184.75
$ws.Range("J122").Value = 6814
$ws.Range("K122").Value = 15554.25
$ws.Range("L122").Value = 20442
$ws.Range("M122").Value = -13104.25
$ws.Range("N122").Value = -25342
$ws.Range("H127").Value = 27777.666
$ws.Range("J127").Value = 27777.666
$ws.Range("L127").Value = 27777.666
$ws.Range("N127").Value = -37697.666

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13749.75
$ws.Range("J15").Value = 39999
$ws.Range("L15").Value = 39999
$ws.Range("N15").Value = -40575
$ws.Range("H41").Value = 13887
$ws.Range("I41").Value = 7775.5
$ws.Range("K41").Value = 7775.5
$ws.Range("M41").Value = -7385.5
$ws.Range("H81").Value = 772688.4
$ws.Range("I81").Value = 2993.5
$ws.Range("J81").Value = 2004200.2
$ws.Range("K81").Value = 5987
$ws.Range("L81").Value = 4008400.4
$ws.Range("M81").Value = -4926
$ws.Range("N81").Value = -4010522.4
$ws.Range("H84").Value = 772688.4
$ws.Range("I84").Value = 2993.5
$ws.Range("J84").Value = 2004200.2
$ws.Range("K84").Value = 29935
$ws.Range("L84").Value = 20042002
$ws.Range("M84").Value = -24631
$ws.Range("N84").Value = -20052610
$ws.Range("H117").Value = 69191.8
$ws.Range("J117").Value = 69191.8
$ws.Range("L117").Value = 69191.8
$ws.Range("N117").Value = -78369.8
$ws.Range("H132").Value = 3250.3333
$ws.Range("I132").Value = 3250.3333
$ws.Range("K132").Value = 9750.999899999999
$ws.Range("M132").Value = -7220.999899999999
